$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B and C headers in rows 2-16 are renamed from "Last year ..." / "This year ..."
# style wording to "Previous ..." / "Upcoming ..." wording (per commit message:
# "change language from lasat year to pervious rotation").

$ws.Range("B2").Value = "Previous cropping rotation"
$ws.Range("C2").Value = "Previous crop"

$ws.Range("B3").Value = "Previous cropping rotation"
$ws.Range("C3").Value = "Previous crop"

$ws.Range("B4").Value = "Previous cropping rotation"
$ws.Range("C4").Value = "Previous crop"

$ws.Range("B5").Value = "Previous cropping rotation"
$ws.Range("C5").Value = "Previous autumn"

$ws.Range("B6").Value = "Previous cropping rotation"
$ws.Range("C6").Value = "Previous autumn"

$ws.Range("B7").Value = "Previous cropping rotation"
$ws.Range("C7").Value = "Previous spring"

$ws.Range("B8").Value = "Upcoming rotation"
$ws.Range("C8").Value = "Upcoming autumn"

$ws.Range("B9").Value = "Previous pasture rotation"
$ws.Range("C9").Value = "Previous pasture rotation"

$ws.Range("B10").Value = "Previous pasture rotation"
$ws.Range("C10").Value = "Previous pasture rotation"

$ws.Range("B11").Value = "Previous pasture rotation"
$ws.Range("C11").Value = "Previous pasture rotation"

$ws.Range("B12").Value = "Previous pasture rotation"
$ws.Range("C12").Value = "Previous autumn"

$ws.Range("B13").Value = "Previous pasture rotation"
$ws.Range("C13").Value = "Previous autumn"

$ws.Range("B14").Value = "Previous pasture rotation"
$ws.Range("C14").Value = "Previous spring"

$ws.Range("B15").Value = "Previous pasture rotation"
$ws.Range("C15").Value = "Previous spring"

$ws.Range("B16").Value = "Previous pasture rotation"
$ws.Range("C16").Value = "Upcoming autumn"

$wb.Save()
